$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that were dropped from the "Spinglass" series
# (A=14 and A=15 in the original data), shifting Girvan Newman and
# Belief rows up by two.
$ws.Rows("16:17").Delete()

# Update the remaining Spinglass values (rows 2-15) with the re-computed numbers
$ws.Range("C3").Value = 7.233346144695181
$ws.Range("C4").Value = 8.896653184233484
$ws.Range("C5").Value = 9.259665523090266
$ws.Range("C6").Value = 13.8051514927761
$ws.Range("C7").Value = 10.78656881226475
$ws.Range("C8").Value = 15.90792291220557
$ws.Range("C9").Value = 10.78203664049489
$ws.Range("C11").Value = 19.16575461971608
$ws.Range("C12").Value = 28.28075184392101
$ws.Range("A14").Value = 12
$ws.Range("C14").Value = 37.11848679514632
$ws.Range("A15").Value = 13
$ws.Range("C15").Value = 74.23697359029265

# Update the Belief series values (rows 22-28) with the re-computed numbers
$ws.Range("C22").Value = 5.164701731932036
$ws.Range("C23").Value = 5.798982448804004
$ws.Range("C24").Value = 6.40475850582917
$ws.Range("C26").Value = 10.08795111505719
$ws.Range("C27").Value = 5.852013683919465
$ws.Range("C28").Value = 9.259665523090266
